$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("D2").Value = 0.287
$ws.Range("G2").Value = 0.1008462623413258
$ws.Range("H2").Value = 0.1008462623413258
$ws.Range("I2").Value = -0.01593794076163611
$ws.Range("J2").Value = -0.01593794076163611
$ws.Range("K2").Value = -23.7
$ws.Range("L2").Value = -0.3342736248236953
$ws.Range("U2").Value = 43.3
$ws.Range("V2").Value = 0.2336751214247166
$ws.Range("W2").Value = -0.08485499462943072
$ws.Range("X2").Value = 0.08643479325601124
$ws.Range("Y2").Value = -0.171289787885442
$ws.Range("Z2").Value = 0.316235504014273
$ws.Range("AA2").Value = -0.00504014272970562
$ws.Range("AB2").Value = 0.08013482295821496
$ws.Range("AC2").Value = -0.08517496568792057
$ws.Range("AD2").Value = 52
$ws.Range("AF2").Value = 52
$ws.Range("AG2").Value = 8.700000000000003
$ws.Range("AH2").Value = 0.2191319005478297
$ws.Range("AI2").Value = 0.16
$ws.Range("AJ2").Value = 0.04484536082474228
$ws.Range("AK2").Value = 0.03088391906283281
$ws.Range("AL2").Value = 8.69
$ws.Range("AM2").Value = 3.449999999999999
$ws.Range("AN2").Value = 9.506398537477148
$ws.Range("AO2").Value = -0.1300345224395857
$ws.Range("AP2").Value = 1.590493601462523
$ws.Range("AQ2").Value = -0.327536231884058

# Row 3 updates
$ws.Range("D3").Value = 0.287
$ws.Range("G3").Value = 0.1008462623413258
$ws.Range("H3").Value = 0.1008462623413258
$ws.Range("I3").Value = -0.01593794076163611
$ws.Range("J3").Value = -0.01593794076163611
$ws.Range("K3").Value = -23.7
$ws.Range("L3").Value = -0.3342736248236953
$ws.Range("U3").Value = 43.3
$ws.Range("V3").Value = 0.2336751214247166
$ws.Range("W3").Value = -0.08485499462943072
$ws.Range("X3").Value = 0.08643479325601124
$ws.Range("Y3").Value = -0.171289787885442
$ws.Range("Z3").Value = 0.316235504014273
$ws.Range("AA3").Value = -0.00504014272970562
$ws.Range("AB3").Value = 0.08013482295821496
$ws.Range("AC3").Value = -0.08517496568792057
$ws.Range("AD3").Value = 52
$ws.Range("AF3").Value = 52
$ws.Range("AG3").Value = 8.700000000000003
$ws.Range("AH3").Value = 0.2191319005478297
$ws.Range("AI3").Value = 0.16
$ws.Range("AJ3").Value = 0.04484536082474228
$ws.Range("AK3").Value = 0.03088391906283281
$ws.Range("AL3").Value = 8.69
$ws.Range("AM3").Value = 3.449999999999999
$ws.Range("AN3").Value = 9.506398537477148
$ws.Range("AO3").Value = -0.1300345224395857
$ws.Range("AP3").Value = 1.590493601462523
$ws.Range("AQ3").Value = -0.327536231884058
